# fix: change not used input file to none.
#
# On sheet "CONDUCTOR_files" the "Value" column (E) for several rows
# referenced placeholder/dummy input file names that are not actually used
# by the code. These are replaced with the literal string "none" to make
# clear that no external file is used for that input.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CONDUCTOR_files")

$ws.Range("E8").Value  = "none"  # EXTERNAL_ALPHAB  (was alphab_dummy.xlsx)
$ws.Range("E9").Value  = "none"  # EXTERNAL_BFIELD  (was bfield.xlsx)
$ws.Range("E10").Value = "none"  # EXTERNAL_CURRENT (was I_file_dummy.xlsx)
$ws.Range("E11").Value = "none"  # EXTERNAL_FLOW    (was flow_dummy.xlsx)
$ws.Range("E12").Value = "none"  # EXTERNAL_HEAT    (was Q_file_dummy.xlsx)
$ws.Range("E13").Value = "none"  # EXTERNAL_STRAIN  (was strain_dummy.xlsx)
$ws.Range("E15").Value = "none"  # EXTERNAL_GRID    (was spatial_discretization.xlsx)

# Match the author's final cursor position on the CONDUCTOR_files sheet.
$ws.Activate() | Out-Null
$ws.Range("G14").Select() | Out-Null
